$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell while forcing text storage so that
# number-like strings (e.g. "73.035.19", "1.00", "0.0490") are not
# auto-converted to numeric values by Excel, matching the source data.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "73.035.19"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.50%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.994.72"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.74%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "619.72"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +15.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.94"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +10.50%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.22%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.761"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.168"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.59"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +8.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000318"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.28"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.47%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.628.31"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.993.75"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.61%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +6.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.36"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.73"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.38%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.908.26"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "442.26"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.95"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +17.40%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "96.68"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.64%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.42"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -5.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.64"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.77%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.17"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.41"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.63"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.50%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.28"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.97%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.84"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -5.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.94"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.70%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "49.26"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.89%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "71.91"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +6.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "643.47"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0911"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +10.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.438"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.23%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.25%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.34"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.94%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.91"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.92%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0490"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.87%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.41"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.926.19"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.55%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.10"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.08%  "

# Rows 47 and 48 swap coin identities with new values
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.97"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +36.65%  "

$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.65"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.03%  "

